$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DeviceId value in A2 (was "10.255.191.231:5556", now "8XV7N16A14001706")
$ws.Range("A2").Value = "8XV7N16A14001706"

# Remove row 3's DeviceId/AndroidPackage values (A3, B3) - row 3 now only keeps K3/L3
$ws.Range("A3:B3").ClearContents()

# Update view: remove topLeftCell freeze on B1, change selection
$ws.Range("A2:B2").Select()
